# Update "想去人数" (want-to-go count) and one ticket-price figure across
# the 展览 / 本地生活 / 全部类型 sheets, matching the refreshed data pull.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 21397
$ws1.Range("F3").Value = 3276
$ws1.Range("F4").Value = 852
$ws1.Range("F6").Value = 546
$ws1.Range("F7").Value = 796
$ws1.Range("F11").Value = 137
$ws1.Range("F12").Value = 562
$ws1.Range("F14").Value = 336
$ws1.Range("G16").Value = 23.3
$ws1.Range("F17").Value = 167
$ws1.Range("F18").Value = 42
$ws1.Range("F20").Value = 77

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 723
$ws3.Range("F4").Value = 725
$ws3.Range("F5").Value = 1725
$ws3.Range("F6").Value = 79

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 723
$ws4.Range("F4").Value = 725
$ws4.Range("F5").Value = 1725
$ws4.Range("F6").Value = 21397
$ws4.Range("F7").Value = 3276
$ws4.Range("F8").Value = 852
$ws4.Range("F10").Value = 79
$ws4.Range("F12").Value = 546
$ws4.Range("F13").Value = 796
$ws4.Range("F20").Value = 137
$ws4.Range("F23").Value = 562
$ws4.Range("F27").Value = 336
$ws4.Range("G30").Value = 23.3
$ws4.Range("F32").Value = 167
$ws4.Range("F33").Value = 42
$ws4.Range("F37").Value = 77
